$d = $word.ActiveDocument

# --- 1. Title line: "Homework #1:" -> "Homework #2" + _GoBack bookmark + ":" ---
# The title run currently reads "Homework #1:" (chars 0-11), followed by a
# " " run and a "Sensors and C++" run. We want the final state to read
# "Homework #2:" again but split as: "Homework #2" | bookmark _GoBack | ":" ,
# with the trailing " " / "Sensors and C++" runs left intact (as separate
# runs, same as before the edit). To stop the engine from silently
# re-merging every same-formatted run in the paragraph when we touch the
# "1" character, we first drop temporary bookmarks right on the existing
# run boundaries (pos 12 and pos 13) so those boundaries survive the edit,
# then do the single-character replacement, add the real _GoBack bookmark
# at the new split point, and finally remove the temporary markers again.

$sep1 = $d.Range(12, 12)
$d.Bookmarks.Add("TempSep1", $sep1)
$sep2 = $d.Range(13, 13)
$d.Bookmarks.Add("TempSep2", $sep2)

$digit = $d.Range(10, 11)
$digit.Text = "2"

$goBackRange = $d.Range(11, 11)
$d.Bookmarks.Add("_GoBack", $goBackRange)

$d.Bookmarks.Item("TempSep1").Delete()
$d.Bookmarks.Item("TempSep2").Delete()

# Note: "_GoBack" is Word's single last-edit-location bookmark, so re-adding
# it above already relocated it away from the end of the document (where it
# used to sit after "...acts gracefully when it is rotated.") to its new
# spot in the title line - nothing further to delete there.
